# Add a new "2021" column (R) to the worksheet, mirroring the existing
# Q (2020) column's layout/formatting, and fill in the new data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new year label
$ws.Range("R4").Value = 2021

# Data rows for the new 2021 column
$ws.Range("R5").Value  = 1
$ws.Range("R6").Value  = 2.2000000000000002
$ws.Range("R7").Value  = 1.7
$ws.Range("R8").Value  = "-"
$ws.Range("R9").Value  = 0.3
$ws.Range("R10").Value = 1.1000000000000001
$ws.Range("R11").Value = "-"
$ws.Range("R12").Value = 0.9
$ws.Range("R13").Value = 0.4
$ws.Range("R14").Value = 0.6

# Copy formatting from column Q so the new column matches the rest of the
# table's visual style (borders, number format, etc.)
$ws.Range("Q4:Q14").Copy()
$ws.Range("R4:R14").PasteSpecial(-4122)  # xlPasteFormats

# Move the active selection similarly to how Excel nudges it after edits
$ws.Range("S17").Select()
